# payments.xlsx — add payment 71652621 (Cash) 2025-08-15T09:35:01
#
# The diff also shows row 5's phone number (A5) switching from a
# text-typed "71652621" to a genuine numeric 71652621, so fix that up
# too before appending the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A5: was stored as text, should be a real number ---
$ws.Range("A5").Value = 71652621

# --- Append new row 6 ---
# A6 keeps the phone number as text (matches the pre-fix formatting that
# row 5 used to have), the rest are plain values.
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "71652621"
$ws.Range("A6").ClearFormats()

$ws.Range("B6").Value = 71
$ws.Range("C6").Value = "Cash"
$ws.Range("D6").Value = "2025-08-15T09:35:01"
